{"js": "const replacements = [\n  [\"652\u00d77=\", \"793\u00d72=\"],\n  [\"776\u00d79=\", \"922\u00d75=\"],\n  [\"419\u00d79=\", \"977\u00d77=\"],\n  [\"633\u00d76=\", \"585\u00d75=\"],\n  [\"186\u00d77=\", \"199\u00d75=\"],\n  [\"377\u00d78=\", \"969\u00d73=\"],\n  [\"884\u00d79=\", \"359\u00d79=\"],\n  [\"864\u00d73=\", \"501\u00d78=\"],\n  [\"571\u00d72=\", \"635\u00d76=\"],\n  [\"273\u00d72=\", \"458\u00d77=\"],\n  [\"739\u00d73=\", \"580\u00d78=\"],\n  [\"964\u00d74=\", \"336\u00d79=\"],\n  [\"286\u00d76=\", \"893\u00d74=\"],\n  [\"960\u00d78=\", \"557\u00d74=\"],\n  [\"891\u00d74=\", \"551\u00d77=\"],\n  [\"899\u00d74=\", \"936\u00d76=\"],\n  [\"104\u00d73=\", \"406\u00d75=\"],\n  [\"146\u00d73=\", \"207\u00d74=\"],\n  [\"137\u00d78=\", \"979\u00d76=\"],\n  [\"943\u00d76=\", \"406\u00d74=\"],\n  [\"554\u00d79=\", \"904\u00d73=\"],\n  [\"281\u00d75=\", \"690\u00d74=\"],\n  [\"473\u00d74=\", \"416\u00d77=\"],\n  [\"259\u00d74=\", \"334\u00d72=\"],\n  [\"423\u00d79=\", \"149\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"652\u00d77=\", \"793\u00d72=\"),\n  @(\"776\u00d79=\", \"922\u00d75=\"),\n  @(\"419\u00d79=\", \"977\u00d77=\"),\n  @(\"633\u00d76=\", \"585\u00d75=\"),\n  @(\"186\u00d77=\", \"199\u00d75=\"),\n  @(\"377\u00d78=\", \"969\u00d73=\"),\n  @(\"884\u00d79=\", \"359\u00d79=\"),\n  @(\"864\u00d73=\", \"501\u00d78=\"),\n  @(\"571\u00d72=\", \"635\u00d76=\"),\n  @(\"273\u00d72=\", \"458\u00d77=\"),\n  @(\"739\u00d73=\", \"580\u00d78=\"),\n  @(\"964\u00d74=\", \"336\u00d79=\"),\n  @(\"286\u00d76=\", \"893\u00d74=\"),\n  @(\"960\u00d78=\", \"557\u00d74=\"),\n  @(\"891\u00d74=\", \"551\u00d77=\"),\n  @(\"899\u00d74=\", \"936\u00d76=\"),\n  @(\"104\u00d73=\", \"406\u00d75=\"),\n  @(\"146\u00d73=\", \"207\u00d74=\"),\n  @(\"137\u00d78=\", \"979\u00d76=\"),\n  @(\"943\u00d76=\", \"406\u00d74=\"),\n  @(\"554\u00d79=\", \"904\u00d73=\"),\n  @(\"281\u00d75=\", \"690\u00d74=\"),\n  @(\"473\u00d74=\", \"416\u00d77=\"),\n  @(\"259\u00d74=\", \"334\u00d72=\"),\n  @(\"423\u00d79=\", \"149\u00d77=\"),\n)\n\nforeach ($p in $pairs) {\n  $oldText = $p[0]\n  $newText = $p[1]\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
